# "Generate Report for Archive"
#
# 1. The localization status moves from "Ready for handoff" to
#    "In Translation" for both target-language columns on the Overview
#    sheet, and for the per-language Status column on the zh-cn / de-de
#    detail sheets.
# 2. The zh-cn / de-de status columns (Overview!E:F, zh-cn!C, de-de!C) are
#    narrowed to match the new, shorter status text.

$wb = $excel.ActiveWorkbook

# --- 1. Update status text -------------------------------------------------

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2. Narrow the status columns ------------------------------------------

$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511

$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511

$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
